$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24: B before A (matches original author's entry order)
$ws.Cells.Item(24, 2).Value = "sqlite3RowSetInsert(RowSet *p, i64 rowid)"
$ws.Cells.Item(24, 1).Value = "pkg1/sqlite-3.21-buggy/src/rowset.c"
$ws.Cells.Item(24, 3).Value = 280

# Row 25: A then B
$ws.Cells.Item(25, 1).Value = "pkg1/sqlite-3.21-buggy/src/tokenize.c"
$ws.Cells.Item(25, 2).Value = "sqlite3IsIdChar(u8 c)"
$ws.Cells.Item(25, 3).Value = 191

# Row 26: B before A
$ws.Cells.Item(26, 2).Value = "sqlite3DropTriggerPtr(Parse *pParse, Trigger *pTrigger)"
$ws.Cells.Item(26, 1).Value = "pkg1/sqlite-3.21-buggy/src/trigger.c"
$ws.Cells.Item(26, 3).Value = 570

# Row 27: A then B
$ws.Cells.Item(27, 1).Value = "pkg1/sqlite-3.21-buggy/src/vdbeapi.c"
$ws.Cells.Item(27, 2).Value = "sqlite3_result_error_code(sqlite3_context *pCtx, int errCode)"
$ws.Cells.Item(27, 3).Value = 490

# Row 28
$ws.Cells.Item(28, 1).Value = "pkg1/sqlite-3.21-buggy/src/vdbeapi.c"
$ws.Cells.Item(28, 2).Value = "sqlite3_bind_zeroblob(sqlite3_stmt *pStmt, int i, int n)"
$ws.Cells.Item(28, 3).Value = 1500

# Row 29
$ws.Cells.Item(29, 1).Value = "pkg1/sqlite-3.21-buggy/src/vdbeaux.c"
$ws.Cells.Item(29, 2).Value = "sqlite3VdbeReusable(Vdbe *p)"
$ws.Cells.Item(29, 3).Value = 413

# Row 30
$ws.Cells.Item(30, 1).Value = "pkg1/sqlite-3.21-buggy/src/vdbesort.c"
$ws.Cells.Item(30, 2).Value = "vdbePmaReaderSeek(SortSubtask *pTask, PmaReader *pReadr, SorterFile *pFile, i64 iOff)"
$ws.Cells.Item(30, 3).Value = 649

# Row 31
$ws.Cells.Item(31, 1).Value = "pkg1/sqlite-3.21-buggy/src/vdbesort.c"
$ws.Cells.Item(31, 2).Value = "sqlite3WalLimit(Wal *pWal, i64 iLimit)"
$ws.Cells.Item(31, 3).Value = 1359

# Row 32
$ws.Cells.Item(32, 1).Value = "pkg1/sqlite-3.21-buggy/src/where.c"
$ws.Cells.Item(32, 2).Value = "sqlite3WhereIsSorted(WhereInfo *pWInfo)"
$ws.Cells.Item(32, 3).Value = 3789

$ws.Application.ActiveWindow.Zoom = 75
$ws.Range("B1").Select()

Write-Host "done"
